# "add suavização de laplace"
# A new row is inserted at the very top of the sheet (pushing the
# existing 200 tweets down by one row) and the word "Teste" is typed
# into the new A1 cell, in bold. Column A is re-selected / re-sized to
# fit its (unchanged) content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 1, shifting everything down.
$ws.Rows.Item(1).Insert()

# Type the new value into the freshly inserted A1 cell and make it bold.
$ws.Range("A1").Value = "Teste"
$ws.Range("A1").Font.Bold = $true

# Select the full column and fit its width to the (longest) content,
# mirroring the column-width/selection state recorded after the edit.
[void]$ws.Columns.Item(1).Select()
[void]$ws.Columns.Item(1).AutoFit()
